$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Sending cluster: ECs)
$ws.Range("G2").Value = 72.92148999999999
$ws.Range("H2").Value = 218.76447
$ws.Range("I2").Value = 0.2015977907456805
$ws.Range("J2").Value = 0.2015977907456805
$ws.Range("M2").Value = 2.546605
$ws.Range("N2").Value = 7.639815
$ws.Range("Q2").Value = 185.70223104145
$ws.Range("R2").Value = 1671.32007937305
$ws.Range("S2").Value = 0.2015977907456805
$ws.Range("T2").Value = 0.2015977907456805

# Row 3 (Sending cluster: FAPs)
$ws.Range("I3").Value = 0.591090693015494
$ws.Range("J3").Value = 0.591090693015494
$ws.Range("M3").Value = 2.546605
$ws.Range("N3").Value = 7.639815
$ws.Range("Q3").Value = 544.4844411975083
$ws.Range("R3").Value = 4900.359970777575
$ws.Range("S3").Value = 0.591090693015494
$ws.Range("T3").Value = 0.591090693015494

# Row 4 (Sending cluster: MuSCs)
$ws.Range("G4").Value = 74.98824566666667
$ws.Range("H4").Value = 224.964737
$ws.Range("I4").Value = 0.2073115162388255
$ws.Range("J4").Value = 0.2073115162388255
$ws.Range("M4").Value = 2.546605
$ws.Range("N4").Value = 7.639815
$ws.Range("Q4").Value = 190.9654413559617
$ws.Range("R4").Value = 1718.688972203655
$ws.Range("S4").Value = 0.2073115162388255
$ws.Range("T4").Value = 0.2073115162388255
